$d = $word.ActiveDocument

function RepeatChar($ch, $n) {
    $s = ""
    for ($i = 0; $i -lt $n; $i++) { $s = $s + $ch }
    return $s
}

# Re-writes the text of a Range in place (forces a real diff so
# adjacent proofErr/grammar markers inside the span are dropped)
# by first swapping in a same-length placeholder, then the real text.
function RewriteRange($theRng, $newText) {
    $n = $newText.Length
    $theRng.Text = RepeatChar "Q" $n
    $theRng.Text = $newText
}

# ---------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently wraps "EEOB590"
# ---------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

Write-Host "Step1 done. EEOB590 para: [$($d.Paragraphs.Item(3).Range.Text)]"

# ---------------------------------------------------------------
# 2) Split the "a data set" run into "a " + "data set" and put a
#    (collapsed) "_GoBack" bookmark at the split point.
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("a data set") | Out-Null
$splitAt = $rng.Start + 2   # right after "a "
$bmRng = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Host "Step2 done. para7: [$($d.Paragraphs.Item(7).Range.Text)]"

# ---------------------------------------------------------------
# 3) After the "...referred to as "mappings". " sentence, add a
#    new sentence: 'Aesthestics are "something you can see". '
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('referred to as "mappings". ') | Out-Null
$rng.Collapse(0)   # wdCollapseEnd
$rng.InsertAfter('Aesthestics')
$rng.Collapse(0)
$rng.InsertAfter(' are ' + [char]8220 + 'something you can see' + [char]8221 + '. ')

Write-Host "Step3 done. para9: [$($d.Paragraphs.Item(9).Range.Text)]"

# ---------------------------------------------------------------
# 4) "Start with ggplot(data, aes())..." -- drop the gramStart/gramEnd
#    proof marks around "(" and merge "(" + "data, " into one run.
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Start with ggplot(data, aes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor = $rng.Start
$parenPos = $anchor + ("Start with ggplot").Length
$mergeRng = $d.Range($parenPos, $parenPos + "(data, ".Length)
Write-Host "mergeRng before: [$($mergeRng.Text)]"
RewriteRange $mergeRng "(data, "

Write-Host "Step4 done. para15: [$($d.Paragraphs.Item(15).Range.Text)]"

# ---------------------------------------------------------------
# 5) "Some plots visualize a transformation  of the original data
#    set." and "Use a stat  to choose a com" -- drop gramStart/End
#    and merge the surrounding runs into single runs.
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Some plots visualize a transformation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart = $rng.Start
$seg1 = "Some plots visualize a transformation  of the original data set."
$seg1Rng = $d.Range($segStart, $segStart + $seg1.Length)
Write-Host "seg1 before: [$($seg1Rng.Text)]"
RewriteRange $seg1Rng $seg1

$rng2 = $d.Content
$rng2.Find.Execute("Use a stat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seg2Start = $rng2.Start
$seg2 = "Use a stat  to choose a com"
$seg2Rng = $d.Range($seg2Start, $seg2Start + $seg2.Length)
Write-Host "seg2 before: [$($seg2Rng.Text)]"
RewriteRange $seg2Rng $seg2

Write-Host "Step5 done. para20: [$($d.Paragraphs.Item(20).Range.Text)]"
